# Auto-generated script applying scheduled market-price refresh to Fenrir_Profits leve data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 13828009
$ws.Range("I70").Value = 41482220
$ws.Range("J70").Value = 904.05554
$ws.Range("K70").Value = 124446660
$ws.Range("L70").Value = 2712.16662
$ws.Range("M70").Value = -124446390
$ws.Range("N70").Value = -3252.16662

$ws.Range("H73").Value = 13828009
$ws.Range("I73").Value = 41482220
$ws.Range("J73").Value = 904.05554
$ws.Range("K73").Value = 124446660
$ws.Range("L73").Value = 2712.16662
$ws.Range("M73").Value = -124445724
$ws.Range("N73").Value = -4584.16662

$ws.Range("H86").Value = 1823.0476
$ws.Range("I86").Value = 1900.3
$ws.Range("J86").Value = 1752.8182
$ws.Range("K86").Value = 1900.3
$ws.Range("L86").Value = 1752.8182
$ws.Range("M86").Value = -777.3
$ws.Range("N86").Value = -3998.8182

$ws.Range("H87").Value = 33115.5
$ws.Range("J87").Value = 33115.5
$ws.Range("L87").Value = 33115.5
$ws.Range("N87").Value = -35611.5

$ws.Range("H89").Value = 1823.0476
$ws.Range("I89").Value = 1900.3
$ws.Range("J89").Value = 1752.8182
$ws.Range("K89").Value = 9501.5
$ws.Range("L89").Value = 8764.091
$ws.Range("M89").Value = -3885.5
$ws.Range("N89").Value = -19996.091

$ws.Range("H90").Value = 33115.5
$ws.Range("J90").Value = 33115.5
$ws.Range("L90").Value = 99346.5
$ws.Range("N90").Value = -111826.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28545.348
$ws.Range("I32").Value = 11503.9795
$ws.Range("J32").Value = 77664.586
$ws.Range("K32").Value = 11503.9795
$ws.Range("L32").Value = 77664.586
$ws.Range("M32").Value = -11216.9795
$ws.Range("N32").Value = -78238.586

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 50997.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 50997.5
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 50997.5
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -52437.5

$ws.Range("H134").Value = 4512670
$ws.Range("I134").Value = 5962394
$ws.Range("J134").Value = 2418.2222
$ws.Range("K134").Value = 17887182
$ws.Range("L134").Value = 7254.6666
$ws.Range("M134").Value = -17884647
$ws.Range("N134").Value = -12324.6666

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 50997.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 50997.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 50997.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -61197.5

$ws.Range("H137").Value = 21599.666
$ws.Range("J137").Value = 21599.666
$ws.Range("L137").Value = 21599.666
$ws.Range("N137").Value = -31799.666

$ws.Range("H140").Value = 32362.637
$ws.Range("J140").Value = 32362.637
$ws.Range("L140").Value = 32362.637
$ws.Range("N140").Value = -42722.637

$ws.Range("H141").Value = 41583.332
$ws.Range("J141").Value = 41583.332
$ws.Range("L141").Value = 41583.332
$ws.Range("N141").Value = -51943.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3497.7
$ws.Range("I86").Value = 2696.6667
$ws.Range("J86").Value = 3841
$ws.Range("K86").Value = 2696.6667
$ws.Range("L86").Value = 3841
$ws.Range("M86").Value = -1573.6667
$ws.Range("N86").Value = -6087

$ws.Range("H89").Value = 3497.7
$ws.Range("I89").Value = 2696.6667
$ws.Range("J89").Value = 3841
$ws.Range("K89").Value = 13483.3335
$ws.Range("L89").Value = 19205
$ws.Range("M89").Value = -7867.333500000001
$ws.Range("N89").Value = -30437

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 723.6875
$ws.Range("I5").Value = 563.9
$ws.Range("J5").Value = 990
$ws.Range("K5").Value = 1691.7
$ws.Range("L5").Value = 2970
$ws.Range("M5").Value = -1579.7
$ws.Range("N5").Value = -3194

$ws.Range("H122").Value = 840.9394
$ws.Range("I122").Value = 833.3200000000001
$ws.Range("J122").Value = 864.75
$ws.Range("K122").Value = 7499.88
$ws.Range("L122").Value = 7782.75
$ws.Range("M122").Value = -5049.88
$ws.Range("N122").Value = -12682.75

$ws.Range("H135").Value = 723.6875
$ws.Range("I135").Value = 563.9
$ws.Range("J135").Value = 990
$ws.Range("K135").Value = 5075.099999999999
$ws.Range("L135").Value = 8910
$ws.Range("M135").Value = -2540.099999999999
$ws.Range("N135").Value = -13980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2579.7
$ws.Range("J80").Value = 2588.5557
$ws.Range("L80").Value = 2588.5557
$ws.Range("N80").Value = -4584.5557

$ws.Range("H83").Value = 2579.7
$ws.Range("J83").Value = 2588.5557
$ws.Range("L83").Value = 12942.7785
$ws.Range("N83").Value = -22926.7785

$ws.Range("H126").Value = 2291.5454
$ws.Range("I126").Value = 1515.3846
$ws.Range("J126").Value = 3412.6667
$ws.Range("K126").Value = 4546.1538
$ws.Range("L126").Value = 10238.0001
$ws.Range("M126").Value = -2076.1538
$ws.Range("N126").Value = -15178.0001

$ws.Range("H140").Value = 14722.5
$ws.Range("J140").Value = 14722.5
$ws.Range("L140").Value = 14722.5
$ws.Range("N140").Value = -25082.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25005062
$ws.Range("I132").Value = 40005040
$ws.Range("J132").Value = 5099.3335
$ws.Range("K132").Value = 120015120
$ws.Range("L132").Value = 15298.0005
$ws.Range("M132").Value = -120012590
$ws.Range("N132").Value = -20358.0005

$ws.Range("H134").Value = 58973
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 58973
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 58973
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -69113

$ws.Range("H135").Value = 22128.5
$ws.Range("J135").Value = 22128.5
$ws.Range("L135").Value = 22128.5
$ws.Range("N135").Value = -32268.5

$ws.Range("H136").Value = 4258.773
$ws.Range("I136").Value = 8336.625
$ws.Range("J136").Value = 1928.5714
$ws.Range("K136").Value = 25009.875
$ws.Range("L136").Value = 5785.7142
$ws.Range("M136").Value = -22459.875
$ws.Range("N136").Value = -10885.7142

$ws.Range("H138").Value = 60214.5
$ws.Range("J138").Value = 60214.5
$ws.Range("L138").Value = 60214.5
$ws.Range("N138").Value = -70494.5

$ws.Range("H139").Value = 60749.5
$ws.Range("J139").Value = 60749.5
$ws.Range("L139").Value = 60749.5
$ws.Range("N139").Value = -71029.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4125
$ws.Range("J62").Value = 4125
$ws.Range("L62").Value = 4125
$ws.Range("N62").Value = -5373

$ws.Range("H65").Value = 4125
$ws.Range("J65").Value = 4125
$ws.Range("L65").Value = 20625
$ws.Range("N65").Value = -26865

$ws.Range("H81").Value = 7144806
$ws.Range("I81").Value = 14286213
$ws.Range("J81").Value = 3398.5715
$ws.Range("K81").Value = 28572426
$ws.Range("L81").Value = 6797.143
$ws.Range("M81").Value = -28571365
$ws.Range("N81").Value = -8919.143

$ws.Range("H84").Value = 7144806
$ws.Range("I84").Value = 14286213
$ws.Range("J84").Value = 3398.5715
$ws.Range("K84").Value = 142862130
$ws.Range("L84").Value = 33985.715
$ws.Range("M84").Value = -142856826
$ws.Range("N84").Value = -44593.715

$ws.Range("H122").Value = 954.8
$ws.Range("I122").Value = 941.3333
$ws.Range("K122").Value = 2823.9999
$ws.Range("M122").Value = -373.9998999999998
